$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the "before" values for the columns that move (D, M, N, O, P, Q, S)
# across rows 2-9, since the edit is a cyclic re-shuffle of these rows'
# weekly records (same market/product, different week -> different date,
# volume and prices).
$snapshot = @{}
for ($r = 2; $r -le 9; $r++) {
    $snapshot[$r] = @{
        D = $ws.Range("D$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        Q = $ws.Range("Q$r").Value2
        S = $ws.Range("S$r").Value2
    }
}

# Destination row -> source row (source row's old data lands in destination row)
$mapping = @{
    2 = 8
    3 = 2
    4 = 3
    5 = 9
    6 = 7
    7 = 4
    8 = 5
    9 = 6
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $data = $snapshot[$srcRow]
    $ws.Range("D$destRow").Value = $data.D
    $ws.Range("M$destRow").Value = $data.M
    $ws.Range("N$destRow").Value = $data.N
    $ws.Range("O$destRow").Value = $data.O
    $ws.Range("P$destRow").Value = $data.P
    $ws.Range("Q$destRow").Value = $data.Q
    $ws.Range("S$destRow").Value = $data.S
}
